$d = $word.ActiveDocument

# --- 1. Remove the "Date du rapport" paragraph together with the blank
#        paragraph that follows it, paragraph marks included, so the text
#        that used to come after collapses directly behind the preceding
#        paragraph (matches the target diff). ---
$dateParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Date du rapport*") {
        $dateParaIndex = $i
        break
    }
}

if ($dateParaIndex -gt 0) {
    $startPara = $d.Paragraphs.Item($dateParaIndex)
    $endPara = $d.Paragraphs.Item($dateParaIndex + 1)
    $killRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $killRange.Delete()
}

# --- 2. Tweak the closing sentence: "exposition nationale." becomes
#        "exposition à l'échelle nationale." ---
$d.Content.Find.Execute("exposition nationale.", $true, $false, $false, $false, $false, $true, 1, $false, "exposition à l’échelle nationale.", 2)
